$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Bump the Version value (row 3, column B) from 0.1.1 to 0.2.0
$ws.Range("B3").Value = "0.2.0"

# 2. Update the publication Date (row 8, column B)
$ws.Range("B8").Value = "2023-10-19T17:05:12+00:00"

# 3. Insert a new "Jurisdiction" row right after "Contact" (row 10) and before
#    "Description" (row 11), pushing everything below down by one row.
$ws.Rows(11).Insert()

# Match the formatting (border/wrap/alignment) used by the other data rows
# instead of the default style Insert() would otherwise apply.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = "iso:code:3166:FR"
